$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "location"
$ws.Range("B1").Value = "price"

$ws.Range("B1").Select()
